$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for Wins, Losses, Ties (columns AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (e.g. AC1) to the new headers
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill season record values (Wins=91, Losses=71, Ties=0) for data rows 2-45
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = 91   # AD
    $ws.Cells.Item($r, 31).Value = 71   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
